# Add new columns I (I0) and J (IF) to Sheet1, mirroring the existing
# H column (IP) header style and the plain numeric data-cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's formatting (bold/border/centered style) onto the two new
# header cells so they reuse the same style definition instead of Excel
# synthesizing a brand-new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (rows 2-17) ----------------------------------------------
# Each entry: row, I-value, J-value (ordered array avoids relying on
# hashtable enumeration order).
$data = @(
    @(2,  9, 9),
    @(3,  6, 6),
    @(4,  7, 7),
    @(5,  8, 8),
    @(6,  8, 8),
    @(7,  7, 8),
    @(8,  7, 7),
    @(9,  8, 8),
    @(10, 9, 9),
    @(11, 7, 7),
    @(12, 8, 8),
    @(13, 7, 7),
    @(14, 6, 6),
    @(15, 6, 6),
    @(16, 8, 8),
    @(17, 5, 5)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}
